$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of a few class name labels in column A
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A4").Value = "mdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"

# Move the active selection to A8
$ws.Range("A8").Select()
